$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.111.97"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "3.580.13"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'577.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.81%  "
$ws.Range("D6").Value = "'186.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.69%  "
$ws.Range("D7").Value = "3.571.06"
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("D8").Value = "'0.619"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.58%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'0.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("E11").Value = "  -4.68%  "
$ws.Range("D12").Value = "'55.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.66%  "
$ws.Range("D13").Value = "'0.0000305"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").Value = "'9.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.78%  "
$ws.Range("D15").Value = "4.150.07"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").Value = "'19.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("D17").Value = "3.580.01"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").Value = "69.958.94"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "'12.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("D22").Value = "'492.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("D24").Value = "'4.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.66%  "
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("D26").Value = "'95.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.65%  "
$ws.Range("D27").Value = "'11.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("E28").Value = "  -7.13%  "
$ws.Range("D29").Value = "'9.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("D30").Value = "'7.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").Value = "'31.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("D32").Value = "'12.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").Value = "'65.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("E34").Value = "  -6.74%  "
$ws.Range("D35").Value = "'579.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.30%  "
$ws.Range("D36").Value = "'3.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +14.57%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").Value = "'38.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.59%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "0.0₃0792"
$ws.Range("E40").Value = "  -5.08%  "
$ws.Range("D41").Value = "'3.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.44%  "
$ws.Range("D42").Value = "'3.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("E43").Value = "  -9.77%  "
$ws.Range("E44").Value = "  -3.58%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.46%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "3.202.85"
$ws.Range("E46").Value = "  -3.96%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0442"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").Value = "'1.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +35.01%  "
$ws.Range("D49").Value = "'9.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("E51").Value = "  +0.04%  "
